# PSP, PMP, meeting report 10.12 update
# - Fill in the two new PSP time-log entries (10/10, 10/12) on the 이준기 sheet
# - Switch the active/selected sheet from 김수인 (sheet2) to 이준기 (sheet1)
# - Restore the two column widths that were nudged in the re-save

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 이준기
$ws2 = $wb.Worksheets.Item(2)   # 김수인

# --- New log row for 10월 10일 ---
$ws1.Range("A16").Value = "10월 10일"
$ws1.Range("B16").Value = 0.5
$ws1.Range("C16").Value = 0.55208333333333337
$ws1.Range("D16").Value = 0
$ws1.Range("E16").Value = 75
$ws1.Range("F16").Value = "Data input, preprocessing, Modelling process Activity Diagram"

# --- New log row for 10월 12일 ---
$ws1.Range("A17").Value = "10월 12일"
$ws1.Range("B17").Value = 0.41666666666666669
$ws1.Range("C17").Value = 0.58333333333333337
$ws1.Range("D17").Value = 60
$ws1.Range("E17").Value = 240
$ws1.Range("F17").Value = "Activity Diagram 수정 및 비교분석 모델 선정"

# --- Column width tweaks on the 이준기 sheet (stored width = ColumnWidth + 0.8333333333333334) ---
$ws1.Columns.Item(4).ColumnWidth = 10.666666666666666
$ws1.Columns.Item(6).ColumnWidth = 45.666666666666664

# --- Selection / active-sheet bookkeeping ---
# Previously 김수인 (sheet2) was the selected/visible tab with F15 selected;
# now 이준기 (sheet1) is the selected tab with F23 selected.
$ws2.Range("F15").Select()
$ws1.Activate()
$ws1.Range("F23").Select()
